$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 4, pushing the old row 4 (k=3, size=4) down to row 5.
$ws.Rows.Item(4).Insert()

# Row 2: k changes from 3 to 4, size from 14 to 5, and the stats are recomputed.
$ws.Range("A2").Value = 4
$ws.Range("C2").Value = 5
$ws.Range("D2").Value = 0.669
$ws.Range("E2").Value = 1.068
$ws.Range("F2").Value = 0.93
$ws.Range("G2").Value = 0.357
$ws.Range("H2").Value = 2.417

# Row 3: k changes from 3 to 4, avg.silwidth is recomputed; size/other stats stay the same.
$ws.Range("A3").Value = 4
$ws.Range("G3").Value = 0.39

# Row 4 (new row, for the new cluster that appears with k=4).
$ws.Range("A4").Value = 4
$ws.Range("B4").Value = "ward.D2"
$ws.Range("C4").Value = 9
$ws.Range("D4").Value = 0.669
$ws.Range("E4").Value = 0.84
$ws.Range("F4").Value = 0.826
$ws.Range("G4").Value = 0.434
$ws.Range("H4").Value = 1.974
$ws.Range("I4").Value = 4

# Row 5 (previously row 4): k changes from 3 to 4, avg.silwidth is recomputed.
$ws.Range("A5").Value = 4
$ws.Range("G5").Value = 0.498
